$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the existing row 110 (the whole
# block of rows 110-232 shifts down to 112-234). Excel's row insert copies
# formatting (e.g. the date number format in column D) from the row above,
# same as interactively inserting rows via the UI.
$ws.Rows("110:111").Insert()

# New row 110: Ajo / Chino / Primera, $/caja 10 kilos, week of 2022-03-09
$ws.Cells.Item(110, 1).Value = 8
$ws.Cells.Item(110, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(110, 3).Value = "Coquimbo"
$ws.Cells.Item(110, 4).Value = 44629
$ws.Cells.Item(110, 5).Value = 4
$ws.Cells.Item(110, 6).Value = 100112003
$ws.Cells.Item(110, 7).Value = "Ajo"
$ws.Cells.Item(110, 8).Value = "Chino"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 400
$ws.Cells.Item(110, 11).Value = 19000
$ws.Cells.Item(110, 12).Value = 20000
$ws.Cells.Item(110, 13).Value = 19500
$ws.Cells.Item(110, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(110, 15).Value = "China"
$ws.Cells.Item(110, 16).Value = 1950
$ws.Cells.Item(110, 17).Value = 10
$ws.Cells.Item(110, 18).Value = "Hortaliza"

# New row 111: Ajo / Chino / Primera, $/malla 10 kilos, week of 2022-03-09
$ws.Cells.Item(111, 1).Value = 8
$ws.Cells.Item(111, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(111, 3).Value = "Coquimbo"
$ws.Cells.Item(111, 4).Value = 44629
$ws.Cells.Item(111, 5).Value = 4
$ws.Cells.Item(111, 6).Value = 100112003
$ws.Cells.Item(111, 7).Value = "Ajo"
$ws.Cells.Item(111, 8).Value = "Chino"
$ws.Cells.Item(111, 9).Value = "Primera"
$ws.Cells.Item(111, 10).Value = 400
$ws.Cells.Item(111, 11).Value = 20000
$ws.Cells.Item(111, 12).Value = 21000
$ws.Cells.Item(111, 13).Value = 20500
$ws.Cells.Item(111, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(111, 15).Value = "China"
$ws.Cells.Item(111, 16).Value = 2050
$ws.Cells.Item(111, 17).Value = 10
$ws.Cells.Item(111, 18).Value = "Hortaliza"
